$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows 15-18 (3./4. smoke test + inheritance page and 2. dialog) ---
# Shared strings must be introduced in this precise order so the underlying
# sharedStrings table lines up: "DefaultInsurance_SmokeTest_FillPage" first,
# then the four per-page smoke-test names (Vehicle/Insurant/Product/SendQuote),
# then the generic "...SmokeTest_FillPage" value.
$ws.Range("D16").Value = "DefaultInsurance_SmokeTest_FillPage"

$ws.Range("A15").Value = "105_CamperInsurance_001_SmokeTest_FillPageVehicleData"
$ws.Range("A16").Value = "105_CamperInsurance_001_SmokeTest_FillPageInsurantData"
$ws.Range("A17").Value = "105_CamperInsurance_001_SmokeTest_FillPageProductData"
$ws.Range("A18").Value = "105_CamperInsurance_001_SmokeTest_FillPageSendQuote"

$ws.Range("C15").Value = "105_CamperInsurance_001_SmokeTest_FillPage"
$ws.Range("E17").Value = "105_CamperInsurance_001_SmokeTest_FillPage"
$ws.Range("G18").Value = "DefaultInsurance_SmokeTest_FillPage"

$ws.Range("B15").Value = "<SET>"
$ws.Range("B16").Value = "<SET>"
$ws.Range("B17").Value = "<SET>"
$ws.Range("B18").Value = "<SET>"

$ws.Range("H15").Value = "<NOP>"
$ws.Range("H16").Value = "<NOP>"
$ws.Range("H17").Value = "<NOP>"
$ws.Range("H18").Value = "<NOP>"

# A17, D16 and E17 use the same "Text" number-format cell style (cellXfs index 1)
# that is already used on A1/B1/A2/B2.
$ws.Range("D16").NumberFormat = "@"
$ws.Range("A17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"

# --- Column width adjustments (best-fit on A, C, E; manual resize on G) ---
# (Target stored widths are 51.5546875 / 40.77734375 / 40.77734375 / 36.33203125 -
#  the values below are the inputs that this host's column-width quantization
#  rounds closest to those numbers.)
$ws.Columns.Item(1).ColumnWidth = 50.666666666666664
$ws.Columns.Item(3).ColumnWidth = 40
$ws.Columns.Item(5).ColumnWidth = 40
$ws.Columns.Item(7).ColumnWidth = 35.5

# --- Selection moved to A13 ---
$ws.Range("A13").Select()

# --- Move/resize the embedded screenshot picture ---
$shp = $ws.Shapes.Item(1)
$shp.Left = 8763000 / 914400 * 72
$shp.Top = 3688080 / 914400 * 72
$shp.Width = 14752381 / 914400 * 72
$shp.Height = 10742857 / 914400 * 72

# --- Window geometry (best effort; this headless host does not persist
#     desktop window placement back into workbookView) ---
$win = $wb.Windows.Item(1)
$win.Width = 33276
$win.Height = 14664
$win.Left = 1116
$win.Top = 1116
